$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, copying the format of the existing
# header cells (e.g. G1: bold font + border + center/top alignment)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H27 with 1 if the "sum" (column G) is greater than 10, else 0
$lastRow = 27
for ($row = 2; $row -le $lastRow; $row++) {
    $sumValue = $ws.Cells.Item($row, 7).Value()
    if ($sumValue -gt 10) {
        $ws.Cells.Item($row, 8).Value = 1
    } else {
        $ws.Cells.Item($row, 8).Value = 0
    }
}
